$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.114.04"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "1.792.81"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.67%  "
$ws.Range("D5").Value = "226.86"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").Value = "0.555"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").Value = "31.08"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").Value = "46.22"
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").Value = "0.280"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").Value = "0.0661"
$ws.Range("E11").Value = "  -1.99%  "
$ws.Range("D12").Value = "0.0927"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "2.051.40"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "11.24"
$ws.Range("E14").Value = "  +9.63%  "
$ws.Range("D15").Value = "1.797.86"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "0.634"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "34.109.15"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "4.21"
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("D19").Value = "69.58"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "253.32"
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("D21").Value = "0.0₃0744"
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").Value = "10.41"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").Value = "4.26"
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").Value = "158.12"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("D27").Value = "16.59"
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("D28").Value = "0.114"
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("D29").Value = "7.01"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "3.90"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").Value = "0.0517"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").Value = "1.20"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").Value = "3.63"
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("D35").Value = "1.87"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").Value = "1.490.65"
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "0.632"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0187"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "83.44"
$ws.Range("E40").Value = "  -5.84%  "
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "0.905"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "2.06"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "0.0516"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D47").Value = "1.948.05"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "5.70"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").Value = "11.77"
$ws.Range("E50").Value = "  +3.81%  "
$ws.Range("D51").Value = "51.31"
$ws.Range("E51").Value = "  -5.65%  "
